$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value2 = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value2 = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value2 = "ds"
$wsForecast.Range("B1").Value2 = "PO_Forecast"
$wsForecast.Range("C1").Value2 = "yhat_lower"
$wsForecast.Range("D1").Value2 = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$dates = 45018.99999999999, 45046.99999999999, 45053.99999999999, 45060.99999999999, 45067.99999999999, 45074.99999999999, 45081.99999999999, 45088.99999999999, 45095.99999999999, 45102.99999999999
$forecast = 1, 1, 1, 1, 1, 1, 1, 1, 1, 1
$lowers = 0.9999999986348875, 0.9999999987441868, 0.9999999984274675, 0.9999999974346818, 0.9999999947267826, 0.999999991219352, 0.9999999867999275, 0.9999999819194079, 0.9999999768617656, 0.9999999717309501
$uppers = 1.000000001276741, 1.000000001346071, 1.000000001561183, 1.000000002198059, 1.000000003815147, 1.000000007050919, 1.00000001108745, 1.00000001503648, 1.00000001934618, 1.00000002401254

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $wsForecast.Cells.Item($r, 1).Value2 = $dates[$i]
    $wsForecast.Cells.Item($r, 2).Value2 = $forecast[$i]
    $wsForecast.Cells.Item($r, 3).Value2 = $lowers[$i]
    $wsForecast.Cells.Item($r, 4).Value2 = $uppers[$i]
}

$wsForecast.Range("A2:A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
